$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date in column C for rows 2-5 by one day
# (2023-10-08 -> 2023-10-09, i.e. serial 45207 -> 45208).
foreach ($r in 2..5) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value()
    $cell.Value = $current.AddDays(1)
}
